# Apply updates to BD1 worksheet: column A year 2015->2017, and new values in column E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3 = @{ A = 20170926; E = 16 }
    4 = @{ A = 20170927; E = 12 }
    5 = @{ A = 20170928; E = 12 }
    6 = @{ A = 20170929; E = 14 }
    7 = @{ A = 20170930; E = 19 }
    8 = @{ A = 20170931; E = 18 }
    9 = @{ A = 20170932; E = 14 }
    10 = @{ A = 20170933; E = 6 }
    11 = @{ A = 20170934; E = 5 }
    12 = @{ A = 20170935; E = 11 }
    13 = @{ A = 20170936; E = 12 }
    14 = @{ A = 20170937; E = 17 }
    15 = @{ A = 20170938; E = 19 }
    16 = @{ A = 20170939; E = 13 }
    17 = @{ A = 20170940; E = 12 }
    18 = @{ A = 20170941; E = 16 }
    19 = @{ A = 20170942; E = 5 }
    20 = @{ A = 20170943; E = 19 }
    21 = @{ A = 20170944; E = 8 }
    22 = @{ A = 20170945; E = 5 }
    23 = @{ A = 20170946; E = 7 }
    24 = @{ A = 20170947; E = 12 }
    25 = @{ A = 20170948; E = 6 }
    26 = @{ A = 20170949; E = 7 }
    27 = @{ A = 20170950; E = 8 }
    28 = @{ A = 20170951; E = 9 }
    29 = @{ A = 20170952; E = 9 }
    30 = @{ A = 20170953; E = 13 }
    31 = @{ A = 20170954; E = 18 }
    32 = @{ A = 20170955; E = 8 }
    33 = @{ A = 20170956; E = 7 }
    34 = @{ A = 20170957; E = 16 }
    35 = @{ A = 20170958 }
    36 = @{ A = 20170959; E = 11 }
    37 = @{ A = 20170960; E = 18 }
    38 = @{ A = 20170961; E = 5 }
    39 = @{ A = 20170962; E = 10 }
    40 = @{ A = 20170963; E = 9 }
    41 = @{ A = 20170964; E = 15 }
    42 = @{ A = 20170965; E = 11 }
    43 = @{ A = 20170966; E = 17 }
    44 = @{ A = 20170967; E = 17 }
    45 = @{ A = 20170968; E = 15 }
    46 = @{ A = 20170969; E = 18 }
    47 = @{ A = 20170970 }
    48 = @{ A = 20170971; E = 9 }
    49 = @{ A = 20170972; E = 15 }
    50 = @{ A = 20170973; E = 13 }
    51 = @{ A = 20170974; E = 8 }
    52 = @{ A = 20170975; E = 16 }
    53 = @{ A = 20170976; E = 16 }
    54 = @{ A = 20170977; E = 6 }
    55 = @{ A = 20170978; E = 16 }
    56 = @{ A = 20170979; E = 7 }
    57 = @{ A = 20170980; E = 11 }
    58 = @{ A = 20170981 }
    59 = @{ A = 20170982; E = 12 }
    60 = @{ A = 20170983; E = 9 }
    61 = @{ A = 20170984; E = 8 }
    62 = @{ A = 20170985; E = 8 }
    63 = @{ A = 20170986; E = 15 }
}

foreach ($row in $changes.Keys) {
    $entry = $changes[$row]
    $ws.Cells.Item([int]$row, 1).Value = $entry.A
    if ($entry.ContainsKey('E')) {
        $ws.Cells.Item([int]$row, 5).Value = $entry.E
    }
}
